# Week 15 logged + Week 16 simulated — appends new-week play-by-play data
# and updates season-to-date aggregate totals across sheets.

$wb = $excel.ActiveWorkbook

# --- YDS sheet: append this week's rush/pass yardage list to each side's string ---
$ws = $wb.Worksheets.Item("YDS")
$ws.Range("B2").Value = $ws.Range("B2").Value2 + " 6 8 11 1 2 2 3 1 5 3 5 7 4 1 -5 11 0 4 6 6 20 -2 2 1"
$ws.Range("C2").Value = $ws.Range("C2").Value2 + " 3 4 0 0 4 10 1 1 12 4 5 0 8 2 5 4 1 6 6 9 3 7 2 0 11 3 2 0 3 24 8 30 4 12 1 4"
$ws.Range("B3").Value = $ws.Range("B3").Value2 + " 16 2 14 19 4 11 6 23 1 9 24 10"
$ws.Range("C3").Value = $ws.Range("C3").Value2 + " 5 1 9 37 15 12 4 27 20 20 7 1 17 9 11 1"

# --- OFF sheet: updated season totals (Home row 2, Road row 3) ---
$ws = $wb.Worksheets.Item("OFF")
$ws.Range("C2").Value = 149
$ws.Range("D2").Value = 9
$ws.Range("F2").Value = 33
$ws.Range("G2").Value = 53
$ws.Range("H2").Value = 6
$ws.Range("J2").Value = 20
$ws.Range("N2").Value = 27
$ws.Range("O2").Value = 22
$ws.Range("B3").Value = 12
$ws.Range("C3").Value = 175
$ws.Range("E3").Value = 33
$ws.Range("F3").Value = 109
$ws.Range("G3").Value = 36
$ws.Range("H3").Value = 23
$ws.Range("I3").Value = 58
$ws.Range("L3").Value = 247
$ws.Range("M3").Value = 145
$ws.Range("Q3").Value = 417

# --- DEF sheet: updated season totals (Home row 2, Road row 3) ---
$ws = $wb.Worksheets.Item("DEF")
$ws.Range("C2").Value = 213
$ws.Range("D2").Value = 16
$ws.Range("F2").Value = 62
$ws.Range("G2").Value = 52
$ws.Range("J2").Value = 28
$ws.Range("N2").Value = 11
$ws.Range("O2").Value = 16
$ws.Range("P2").Value = 10
$ws.Range("C3").Value = 128
$ws.Range("D3").Value = 5
$ws.Range("E3").Value = 30
$ws.Range("F3").Value = 89
$ws.Range("I3").Value = 48
$ws.Range("J3").Value = 51
$ws.Range("L3").Value = 230
$ws.Range("M3").Value = 159
$ws.Range("Q3").Value = 454

# --- ST sheet: updated season totals + appended punt/kick distance lists ---
$ws = $wb.Worksheets.Item("ST")
$ws.Range("B2").Value = 53
$ws.Range("D2").Value = 54
$ws.Range("F2").Value = 19
$ws.Range("G2").Value = 16
$ws.Range("J2").Value = 14
$ws.Range("K2").Value = 14
$ws.Range("B3").Value = 35
$ws.Range("D3").Value = $ws.Range("D3").Value2 + " 48 52 50 29 37 63"
$ws.Range("D4").Value = $ws.Range("D4").Value2 + " 0 0 0 0 0 5"
$ws.Range("D5").Value = $ws.Range("D5").Value2 + " 0 11 12"

# --- TURNS sheet: updated season totals (Road row 3) ---
$ws = $wb.Worksheets.Item("TURNS")
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 4
$ws.Range("E3").Value = 11

# --- PEN sheet: updated season totals ---
$ws = $wb.Worksheets.Item("PEN")
$ws.Range("B3").Value = 14
$ws.Range("D4").Value = 7
